# Upload/refresh of the dashboard data: rescale the "Actual Cost to Date"
# column from millions of N$ to full N$ values, widen a few columns that
# were previously left at default width, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-scale column R (Actual Cost to Date (Mil)) from millions to absolute
# values for every data row (rows 2 through 86).
for ($row = 2; $row -le 86; $row++) {
    $cell = $ws.Cells.Item($row, 18)
    $cell.Value = $cell.Value2 * 1000000
}

# Widen columns I, K, L and M which previously used the default width.
$ws.Columns.Item(9).ColumnWidth = 20.166666666666668
$ws.Columns.Item(11).ColumnWidth = 11.666666666666666
$ws.Columns.Item(12).ColumnWidth = 16.0
$ws.Columns.Item(13).ColumnWidth = 17.333333333333332

# Move the current selection to L1.
$ws.Range("L1").Select()
